$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Append a new row (row 22) continuing the schedule:
#   A22 -> "8/22"  (next entry in the recurring weekly date column)
#   B22 -> "10/17" (next entry in the running date column)
#   C22 -> "第72期 第七代坐騎" (new period description)
$ws.Range("A22").Value = "8/22"
$ws.Range("B22").Value = "10/17"
$ws.Range("C22").Value = "第72期 第七代坐騎"

# Match the saved view state: scrolled so row 6 is the top visible row,
# and the active selection moved to C20.
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C20").Select() | Out-Null
